# Fruta / hortaliza, semanal
# Weekly data refresh: a new weekly observation is inserted as row 4
# (shifting the previously existing rows 4-8 down to rows 5-9), and the
# new row 4 is populated with this week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4-8 down to 5-9.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44414
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112017
$ws.Range("G4").Value = "Ramas de apio"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 6500
$ws.Range("N4").Value = "$/atado 7 kilos"
$ws.Range("O4").Value = "Región de Arica y Parinacota"
$ws.Range("P4").Value = 6500
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
